$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text (t="inlineStr" in the
# original workbook). Force the number format to Text before assigning so
# Excel does not silently convert the string into a real number, then restore
# the original style so no stray formatting is left behind.
$dRange = $ws.Range("D2:D50")
$dOrigStyle = $dRange.Style
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '247.63'
$ws.Range("D3").Value = '21.79'
$ws.Range("D4").Value = '5.463'
$ws.Range("D5").Value = '0.05694'
$ws.Range("D6").Value = '3.380'
$ws.Range("D7").Value = '0.8049'
$ws.Range("D8").Value = '1.039'
$ws.Range("D9").Value = '0.1486'
$ws.Range("D10").Value = '0.07312'
$ws.Range("D11").Value = '0.03167'
$ws.Range("D12").Value = '0.02937'
$ws.Range("D13").Value = '0.09286'
$ws.Range("D15").Value = '3.383'
$ws.Range("D16").Value = '0.04720'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = '0.0005858'
$ws.Range("E17").Value = '16OneONE'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = '0.006316'
$ws.Range("E18").Value = '17TigerCashTCH'
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").Value = '0.005049'
$ws.Range("E19").Value = '18HotbitTokenHTBBestin24h'
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").Value = '0.001046'
$ws.Range("E20").Value = '19BitKanKAN'
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").Value = '0.0001500'
$ws.Range("E21").Value = '20NitroExNTX'
$ws.Range("B22").Value = 'UpBots'
$ws.Range("C22").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D22").Value = '0.0003199'
$ws.Range("E22").Value = '21UpBotsUBXT'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = '3.774'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'KuCoinToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D24").Value = '6.427'
$ws.Range("E24").Value = '23KuCoinTokenKCS'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = '2.123'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("D40").Value = '0.04114'
$ws.Range("D41").Value = '0.006950'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '0.1044'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '0.002970'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = '0.008111'
$ws.Range("D45").Value = '0.00005831'
$ws.Range("D47").Value = '0.0005498'
$ws.Range("D48").Value = '0.6822'
$ws.Range("D49").Value = '0.009648'
$ws.Range("D50").Value = '0.00002099'

# Restore the original (default) style on column D now that the text values
# are safely stored, so no extra formatting/style is introduced.
$dRange.Style = $dOrigStyle

Write-Host "Applied symbol list update"
